# Auto-generated edit script: update crypto price/volume values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.345.66"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "3.217.76"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.215.98"
$ws.Range("E9").Value = "  -1.95%  "
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.69"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.84%  "
$ws.Range("E12").Value = "  -3.06%  "
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("D15").Value = "3.745.45"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").Value = "66.439.83"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").Value = "3.215.18"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("E18").Value = "  -3.11%  "
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "506.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.86%  "
$ws.Range("E22").Value = "  -2.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.155"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +70.96%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.89%  "
$ws.Range("E30").Value = "  -1.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  -5.33%  "
$ws.Range("E36").Value = "  -2.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "499.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("D39").Value = "0.0₃0770"
$ws.Range("E39").Value = "  +11.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0417"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.27%  "
$ws.Range("E42").Value = "  +3.87%  "
$ws.Range("E43").Value = "  -2.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.293"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.95%  "
$ws.Range("D45").Value = "2.917.95"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.09%  "
